$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value to a cell without Excel re-interpreting
# numeric-looking strings (e.g. "1.002", "28.526.36") as numbers, and without
# leaving the cells style pointed at a new "Text" number-format record -
# we stamp @ on just long enough to take the string literally, then restore
# the original (default) style so the saved XML has no style attribute,
# exactly like the source cells.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "28.526.36"
$ws.Range("E2").Value = "  +0.74%  "
Set-TextValue $ws.Range("D3") "1.919.64"
$ws.Range("E3").Value = "  +2.17%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "316.20"
$ws.Range("E5").Value = "  +1.24%  "
Set-TextValue $ws.Range("D6") "1.002"
$ws.Range("E6").Value = "  +0.03%  "
Set-TextValue $ws.Range("D7") "0.5114"
$ws.Range("E7").Value = "  +1.24%  "
Set-TextValue $ws.Range("D8") "0.3969"
$ws.Range("E8").Value = "  +0.20%  "
Set-TextValue $ws.Range("D9") "0.09731"
$ws.Range("E9").Value = "  -0.93%  "
Set-TextValue $ws.Range("D10") "1.144"
$ws.Range("E10").Value = "  +1.45%  "
Set-TextValue $ws.Range("D11") "42.09"
$ws.Range("E11").Value = "  +1.20%  "
Set-TextValue $ws.Range("D12") "6.474"
$ws.Range("E12").Value = "  +0.04%  "
Set-TextValue $ws.Range("D13") "21.01"
$ws.Range("E13").Value = "  +0.07%  "
Set-TextValue $ws.Range("D14") "1.915.35"
$ws.Range("E14").Value = "  +2.15%  "
Set-TextValue $ws.Range("D15") "7.413"
Set-TextValue $ws.Range("D16") "1.002"
$ws.Range("E16").Value = "  -0.02%  "
Set-TextValue $ws.Range("D17") "0.00001132"
$ws.Range("E17").Value = "  -0.75%  "
Set-TextValue $ws.Range("D18") "94.07"
$ws.Range("E18").Value = "  +0.45%  "
Set-TextValue $ws.Range("D19") "0.06669"
$ws.Range("E19").Value = "  -0.06%  "
Set-TextValue $ws.Range("D20") "18.13"
$ws.Range("E20").Value = "  +4.02%  "
Set-TextValue $ws.Range("D21") "1.001"
$ws.Range("E21").Value = "  -0.07%  "
Set-TextValue $ws.Range("D22") "6.286"
$ws.Range("E22").Value = "  +2.73%  "
Set-TextValue $ws.Range("D23") "28.594.68"
$ws.Range("E23").Value = "  +0.76%  "
Set-TextValue $ws.Range("D24") "11.48"
$ws.Range("E24").Value = "  +1.17%  "
Set-TextValue $ws.Range("D25") "2.308"
$ws.Range("E25").Value = "  +2.09%  "
Set-TextValue $ws.Range("D26") "2.676"
$ws.Range("E26").Value = "  +4.39%  "
Set-TextValue $ws.Range("D27") "2.139.01"
$ws.Range("E27").Value = "  +2.33%  "
Set-TextValue $ws.Range("D28") "21.22"
$ws.Range("E28").Value = "  -0.99%  "
Set-TextValue $ws.Range("D29") "158.19"
$ws.Range("E29").Value = "  +0.39%  "
Set-TextValue $ws.Range("D30") "128.82"
$ws.Range("E30").Value = "  +0.88%  "
Set-TextValue $ws.Range("D31") "1.108"
$ws.Range("E31").Value = "  +4.24%  "
Set-TextValue $ws.Range("D32") "0.1069"
$ws.Range("E32").Value = "  +0.57%  "
Set-TextValue $ws.Range("D33") "5.702"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("E34").Value = "  +0.80%  "
Set-TextValue $ws.Range("D35") "9.790"
$ws.Range("E35").Value = "  +2.87%  "
Set-TextValue $ws.Range("D36") "0.06709"
$ws.Range("E36").Value = "  -1.74%  "
Set-TextValue $ws.Range("D37") "0.02441"
$ws.Range("E37").Value = "  +2.08%  "
Set-TextValue $ws.Range("D38") "1.259"
$ws.Range("E38").Value = "  +3.81%  "
Set-TextValue $ws.Range("D39") "0.2221"
$ws.Range("E39").Value = "  +1.53%  "
Set-TextValue $ws.Range("D40") "11.66"
$ws.Range("E40").Value = "  +1.12%  "
Set-TextValue $ws.Range("D41") "0.6446"
$ws.Range("E41").Value = "  +1.95%  "
Set-TextValue $ws.Range("D42") "5.057"
$ws.Range("E42").Value = "  +0.49%  "
Set-TextValue $ws.Range("D43") "1.206"
$ws.Range("E43").Value = "  +2.42%  "
Set-TextValue $ws.Range("D44") "1.001"
$ws.Range("E44").Value = "  +0.05%  "
Set-TextValue $ws.Range("D45") "13.71"
$ws.Range("E45").Value = "  +1.51%  "
Set-TextValue $ws.Range("D46") "0.6078"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  +2.66%  "
Set-TextValue $ws.Range("D48") "1.284"
$ws.Range("E48").Value = "  +0.32%  "
Set-TextValue $ws.Range("D49") "2.063"
$ws.Range("E49").Value = "  +3.36%  "
Set-TextValue $ws.Range("D50") "124.11"
$ws.Range("E50").Value = "  -1.08%  "
Set-TextValue $ws.Range("D51") "1.198"
$ws.Range("E51").Value = "  -0.35%  "
